{"js": "// Rewrite the bullet list under \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\"\n// from job-duty statements into impact-focused accomplishment statements,\n// per the commit's \"Fix Key Achievements to use proper accomplishment\n// statements\" rewrite. Three bullets are rewritten in place, two bullets\n// are removed outright, and the remaining bullet is rewritten.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text,items/style\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then restrict all work\n// to the block of paragraphs between it and the next top-level (Heading 2)\n// heading -- several of the bullet strings we need to touch also appear\n// verbatim earlier in the document (e.g. under PROFESSIONAL EXPERIENCE),\n// so scoping by section keeps this from touching the wrong occurrence.\nlet sectionStart = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].style === \"Heading 2\" && items[i].text.trim() === \"KEY ACHIEVEMENTS AND IMPACT\") {\n    sectionStart = i;\n    break;\n  }\n}\nif (sectionStart === -1) {\n  throw new Error('Could not find the \"KEY ACHIEVEMENTS AND IMPACT\" heading.');\n}\n\nlet sectionEnd = items.length;\nfor (let i = sectionStart + 1; i < items.length; i++) {\n  if (items[i].style === \"Heading 2\") {\n    sectionEnd = i;\n    break;\n  }\n}\n\n// Old bullet text -> new bullet text (null means \"delete this paragraph\").\nconst replacements = [\n  [\n    \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over $2 trillion\",\n    \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n  ],\n  [\n    \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\",\n    \"\u2022 $4.7M savings enabled nonprofit access\"\n  ],\n  [\n    \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\",\n    \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n  ],\n  [\n    \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations $4.7M and enabling smaller nonprofits to conduct analysis\",\n    null\n  ],\n  [\n    \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\",\n    null\n  ],\n  [\n    \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\",\n    \"\u2022 178% accuracy improvement in racial classification algorithms\"\n  ]\n];\n\nfor (const [oldText, newText] of replacements) {\n  let target = null;\n  for (let i = sectionStart + 1; i < sectionEnd; i++) {\n    if (items[i].text === oldText) {\n      target = items[i];\n      break;\n    }\n  }\n  if (!target) {\n    throw new Error(\"Could not find expected bullet: \" + oldText);\n  }\n  if (newText === null) {\n    target.delete();\n  } else {\n    target.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Rewrite the bullet list under \"KEY ACHIEVEMENTS AND IMPACT\" > \"Impact\"\n# from job-duty statements into impact-focused accomplishment statements,\n# per the commit's \"Fix Key Achievements to use proper accomplishment\n# statements\" rewrite. Three bullets are rewritten in place, two bullets\n# are removed outright, and the remaining bullet is rewritten.\n\n$d = $word.ActiveDocument\n\n# Locate the \"KEY ACHIEVEMENTS AND IMPACT\" heading, then restrict all work to\n# the block of paragraphs between it and the next top-level (Heading 2)\n# heading -- several of the bullet strings we need to touch also appear\n# verbatim earlier in the document (e.g. under PROFESSIONAL EXPERIENCE), so\n# scoping by section keeps this from touching the wrong occurrence.\n$count = $d.Paragraphs.Count\n\n$sectionStart = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`n\")\n    if (($d.Paragraphs.Item($i).Style.NameLocal -eq \"Heading 2\") -and ($txt -eq \"KEY ACHIEVEMENTS AND IMPACT\")) {\n        $sectionStart = $i\n        break\n    }\n}\nif ($sectionStart -eq -1) {\n    throw 'Could not find the \"KEY ACHIEVEMENTS AND IMPACT\" heading.'\n}\n\n$sectionEnd = $count\nfor ($i = $sectionStart + 1; $i -le $count; $i++) {\n    if ($d.Paragraphs.Item($i).Style.NameLocal -eq \"Heading 2\") {\n        $sectionEnd = $i - 1\n        break\n    }\n}\n\n# Old bullet text -> new bullet text ($null means \"delete this paragraph\").\n$replacements = @(\n    @{\n        Old = \"\u2022 Built real-time FEC analysis systems using Python, Pandas and PySpark to detect likely fraud, money laundering and financial crimes across billions of records daily, performing time series analysis on trillions of records in the political spending sub-economy valued over `$2 trillion\"\n        New = \"\u2022 Algorithmic innovation: Pioneered trigonometric boundary estimation reducing mapping costs 73.5%\"\n    },\n    @{\n        Old = \"\u2022 Built cloud-based data warehouse solutions on AWS processing billions of records with 99.94% accuracy\"\n        New = \"\u2022 `$4.7M savings enabled nonprofit access\"\n    },\n    @{\n        Old = \"\u2022 Designed ETL pipelines using PySpark, dbt, and PostgreSQL/PostGIS for large-scale geospatial datasets\"\n        New = \"\u2022 Breakthrough demographic discovery: Uncovered systematic voter miscoding affecting millions\"\n    },\n    @{\n        Old = \"\u2022 Trigonometric algorithm for boundary estimation reduced mapping costs by 73.5%, saving campaigns and organizations `$4.7M and enabling smaller nonprofits to conduct analysis\"\n        New = $null\n    },\n    @{\n        Old = \"\u2022 Built redistricting platform used by thousands of analysts nationwide with real-time collaborative editing and Census integration, serving 12,847 analysts across 89 organizations\"\n        New = $null\n    },\n    @{\n        Old = \"\u2022 Discovered systematic race coding errors affecting all Black and Asian-American voters, developed geospatial machine learning algorithms improving classification accuracy from 23% to 64%\"\n        New = \"\u2022 178% accuracy improvement in racial classification algorithms\"\n    }\n)\n\nforeach ($r in $replacements) {\n    $targetIndex = -1\n    for ($i = $sectionStart + 1; $i -le $sectionEnd; $i++) {\n        $txt = $d.Paragraphs.Item($i).Range.Text.TrimEnd(\"`r\", \"`n\")\n        if ($txt -eq $r.Old) {\n            $targetIndex = $i\n            break\n        }\n    }\n    if ($targetIndex -eq -1) {\n        throw \"Could not find expected bullet: \" + $r.Old\n    }\n\n    if ($r.New -eq $null) {\n        $d.Paragraphs.Item($targetIndex).Range.Delete()\n        $sectionEnd = $sectionEnd - 1\n    } else {\n        $d.Paragraphs.Item($targetIndex).Range.Text = $r.New\n    }\n}\n"}
